# [PHOENIX-5858] Added some files for test
#
# Adds a new "deputyExecutiveEngineer" login-test-data row (row 16) to the
# registeredUserDetails sheet, mirroring the existing rows: dataName / id /
# password / hasZone (a FALSE() formula), plus the matching mailto hyperlink
# on the password cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registeredUserDetails")

# --- new data row --------------------------------------------------------
# Match the row height used by the row just above it.
$ws.Rows(16).RowHeight = $ws.Rows(15).RowHeight

$ws.Range("A16").Value = "deputyExecutiveEngineer"

# Keep the leading zero - this column stores ids as text everywhere else.
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "0942870"

$ws.Range("C16").Value = "kurnool_eGov@123"

# Same "TRUE/FALSE" custom display format used by the other hasZone cells.
$ws.Range("D16").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("D16").Formula = "=FALSE()"

# --- hyperlink on the password cell, matching rows 11-15 -----------------
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123") | Out-Null
# Adding the hyperlink auto-applies Excel's blue/underlined "Hyperlink"
# look; the other password cells in this sheet keep the plain "Normal"
# look even though they carry a hyperlink, so put it back.
$ws.Range("C16").Style = "Normal"

# --- selection / view, matches the author's final cursor position --------
$ws.Range("B16").Select() | Out-Null
